$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 26 stats for 2026-01 (authorities and users counts revised)
$ws.Range("C26").Value = 1011
$ws.Range("D26").Value = 6044688

# Recompute derived columns to match the new authority/user counts
$ws.Range("E26").Value = 931.672009864365
$ws.Range("G26").Value = 7.32484076433122
$ws.Range("H26").Value = 25.88083124371474
